# Apply the 2023-01-20 bread_coop crawl update:
#  1. Bump ratingAmount (D) / ratingValue (E) for a handful of rows whose
#     review counts ticked up between the two crawl runs.
#  2. Swap the two rows that the crawler re-ordered (same two products,
#     positions exchanged) - ids 4944608 (Schaer) and 5909120 (Rob&Lissy).
#  3. Flag the now out-of-stock "Fine Food Sesame Rice Crackers" listing.
#  4. Refresh the timestamp column (O) on every data row to the new crawl
#     time.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. ratingAmount / ratingValue bumps -----------------------------------
# row -> new ratingAmount (D)
$dUpdates = @{
  5   = 28
  9   = 28
  19  = 37
  26  = 40
  33  = 13
  39  = 16
  46  = 12
  50  = 44
  60  = 7
  63  = 12
  77  = 19
  80  = 38
  89  = 43
  92  = 11
  101 = 12
  111 = 4
  124 = 9
  131 = 10
  135 = 11
  156 = 15
  164 = 4
  169 = 29
  170 = 12
  180 = 15
  270 = 4
  328 = 4
}

foreach ($row in $dUpdates.Keys) {
  $ws.Cells.Item($row, 4).Value = $dUpdates[$row]
}

# row -> new ratingValue (E) - only a few rows also changed their rating
$eUpdates = @{
  101 = 4
  180 = 4.5
  270 = 3.5
}

foreach ($row in $eUpdates.Keys) {
  $ws.Cells.Item($row, 5).Value = $eUpdates[$row]
}

# --- 2. swap rows 121 and 122 ----------------------------------------------
# Columns A-N hold product data (O is the timestamp, refreshed below for
# every row anyway). Price-ish text columns (A, G, J) look numeric, so force
# them to stay text the way the source file stored them.
$textForceCols = @(1, 7, 10)   # A, G, J

function Get-RowValues($row) {
  $vals = @{}
  for ($col = 1; $col -le 14; $col++) {
    $vals[$col] = $ws.Cells.Item($row, $col).Value()
  }
  return $vals
}

function Set-RowValues($row, $vals) {
  for ($col = 1; $col -le 14; $col++) {
    $cell = $ws.Cells.Item($row, $col)
    if ($textForceCols -contains $col) {
      $cell.NumberFormat = "@"
    }
    $cell.Value = $vals[$col]
  }
}

$row121 = Get-RowValues 121
$row122 = Get-RowValues 122

Set-RowValues 121 $row122
Set-RowValues 122 $row121

# --- 3. Fine Food Sesame Rice Crackers now flagged as out of stock --------
$ws.Range("M342").Value = "Fine Food Sesame Rice Crackers - Online kein Bestand 4.95 Schweizer Franken"

# --- 4. refresh crawl timestamp on every data row --------------------------
$newTimestamp = "2023-01-20 12:56:40"
$lastRow = $ws.UsedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
  $ws.Cells.Item($row, 15).Value = $newTimestamp
}

Write-Output "done"
